$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "FT231680GNW5GHWN",
    "FT231680XLHY7PKL",
    "FT231680LNFJP367",
    "FT231680GNW5GP9F",
    "FT2316803Z2HPVN4",
    "FT231680MKQM392C",
    "FT2316803Z2HPZDC",
    "FT2316805D7D4GNY",
    "FT231680T6166DLZ",
    "FT231680MKQM407H",
    "FT231680PGLP2QGQ"
)

$startRow = 14
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
